# Cost Budget Estimate.xlsx update
# - Developer hourly rate (C7) increased from 100 to 200
# - Product Owner hourly rate (C9) increased from 100 to 375
# Downstream formulas (D7, F7, D9, F9, B1 total) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C7").Value = 200
$ws.Range("C9").Value = 375

# Leave the selection where the author left it when they saved the file
$ws.Range("B8").Select()
